$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.066.46"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "1.680.05"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'215.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.254"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").Value = "'21.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.68%  "
$ws.Range("D10").Value = "'0.0624"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.58%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "1.916.53"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "1.648.16"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").Value = "'0.534"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "'66.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "27.056.04"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").Value = "'8.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("D19").Value = "'236.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("D23").Value = "'9.27"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").Value = "'147.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").Value = "'7.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("D27").Value = "'16.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "'1.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").Value = "1.546.04"
$ws.Range("E33").Value = "  +6.21%  "
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").Value = "  +4.91%  "
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D38").Value = "'0.916"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("E40").Value = "  +7.26%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'67.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").Value = "1.821.66"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "'0.779"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "'90.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("E51").Value = "  +7.08%  "
